$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell E1, matching the style of the existing header row (D1)
$ws.Range("E1").Value = "Colocação"
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122) # xlPasteFormats

# Ranking values for rows 2-8 (row 9 and 10 are left untouched)
$ws.Range("E2").Value = "1º"
$ws.Range("E3").Value = "2º"
$ws.Range("E4").Value = "3º"
$ws.Range("E5").Value = "4º"
$ws.Range("E6").Value = "5º"
$ws.Range("E7").Value = "6º"
$ws.Range("E8").Value = "16º"
